$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to be treated as text so numeric-looking
# strings (e.g. "228.56", "0.580") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Updated price / 1h-volume figures
$ws.Range("D2").Value = "34.489.36"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.812.56"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "228.56"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "0.580"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "35.89"
$ws.Range("E8").Value = "  +9.12%  "
$ws.Range("D9").Value = "0.304"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("D10").Value = "0.0698"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("D11").Value = "0.0963"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "2.075.12"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "11.53"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "1.824.12"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "0.649"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "4.54"
$ws.Range("E16").Value = "  +5.73%  "
$ws.Range("D17").Value = "34.475.81"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "69.42"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "247.87"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "0.0₃0802"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "11.58"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "4.23"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").Value = "172.30"
$ws.Range("D25").Value = "2.14"
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("D26").Value = "8.00"
$ws.Range("E26").Value = "  +8.40%  "
$ws.Range("D27").Value = "16.94"
$ws.Range("E27").Value = "  +2.20%  "
$ws.Range("E28").Value = "  +3.58%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "4.09"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").Value = "3.89"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("D32").Value = "0.0535"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "1.398.70"
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").Value = "0.682"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "83.99"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").Value = "0.969"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D45").Value = "13.46"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").Value = "0.0505"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "1.973.90"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "105.91"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "0.0₆0130"
$ws.Range("E51").Value = "  +0.22%  "

# Rows 42-44 reordered: WEMIXToken moves up to row 42, MXToken to row 43,
# HuobiToken to row 44 (rank numbers in column A stay as-is).
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.23"
$ws.Range("E42").Value = "  +11.45%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").Value = "  +1.80%  "

$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  +0.23%  "

# Restore default (unstyled) formatting on the touched range so the
# cells keep their original look (no explicit style / number format).
$ws.Range("D2:E51").Style = "Normal"
